# Applies the commit "Updated docx and md files" to the media-analysis docx:
# rewrites the Short Summary, Summary, Question 1, Question 2, Question 4 and
# Entities body paragraphs with the new ChatGPT/Arab News media-analysis copy.
$d = $word.ActiveDocument
$vt = [char]11   # vertical tab == Word's in-paragraph line break (<w:br/>) when read/written via Range.Text

function Set-ParagraphText($anchorText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find anchor text: $anchorText"
    }
    $rng.Expand(4) | Out-Null   # wdParagraph
    $rng.Text = $newText
}

# Short Summary
Set-ParagraphText 'The article discusses the capabilities and potential impacts of ChatGPT, an advanced AI language model developed by OpenAI.' ('The article discusses the advancements in AI technology, particularly focusing on OpenAI s language model, ChatGPT.')

# Summary
Set-ParagraphText 'The article discusses the evolution and potential applications of AI, particularly focusing on OpenAI''s ChatGPT model.' ('In summary:' + $vt + '' + $vt + '1. ChatGPT describes itself as an advanced AI language model developed by OpenAI, capable of providing a wide range of responses based on the patterns and relationships it has learned from its training data.' + $vt + '2. It is being interviewed by Arab News, a leading English-language daily newspaper in Saudi Arabia, which provides news and analysis on local, regional and international events.' + $vt + '3. Microsoft has integrated ChatGPT into its Bing search engine, Edge browser and other products, sparking a race with Google, Amazon, Baidu and Meta to avoid being left behind in the AI technology field.' + $vt + '4. On March 14, OpenAI released an update called GPT-4, which promises more accurate text responses and the ability to generate responses from both image and text inputs.' + $vt + '5. ChatGPT is not capable of predicting the specific features or capabilities that will be added to it in the future, but areas of focus may include improving the naturalness and sophistication of AI-generated text, developing AI models that are more context-aware, and integrating AI models with other data sources for personalization.' + $vt + '6. ChatGPT emphasizes that the impact of AI on society will depend on how it is used by developers, policymakers, and users, and expresses cautious optimism about the potential positive impacts of AI, such as improved healthcare, education, and communication, while also acknowledging concerns about privacy, bias, and job displacement.' + $vt + '7. ChatGPT encourages responsible and ethical use of AI technology for a better future for all.')

# Question 1
Set-ParagraphText 'In the sample analyzed, the media present ChatGPT as a tool with significant potential benefits and concerns.' ('The media''s framing of public discussions about ChatGPT tends to emphasize its technological advancements, potential benefits, and the ongoing competition among tech giants in the AI sector. In some instances, it is depicted as a game-changer or disruptor, while also acknowledging concerns about privacy, bias, and job displacement.' + $vt + '' + $vt + 'As for metaphors that keep cropping up, here are a few examples:' + $vt + '' + $vt + '1. **Revolutionary Force**: ChatGPT and AI in general are frequently described as revolutionizing various industries and transforming the way people live, work, and interact with technology. This metaphor highlights the profound impact that AI may have on society.' + $vt + '2. **Race/Competition**: The development of AI is often framed as a competition or race among tech companies like Google, Microsoft, Amazon, Baidu, and Meta. This metaphor emphasizes the urgency and importance of being at the forefront of AI technology.' + $vt + '3. **Ticking Time Bomb**: Some reports discuss concerns about the potential negative consequences of AI, likening it to a ticking time bomb that could cause significant harm if not properly managed or regulated. This metaphor underscores the need for responsible development and deployment of AI technologies.' + $vt + '4. **Double-edged Sword**: The dual nature of AI is frequently described as a double-edged sword, emphasizing both its potential benefits and risks. This metaphor acknowledges that while AI has tremendous promise, it also carries inherent dangers that must be carefully addressed.' + $vt + '5. **Human-like/Superintelligent Being**: ChatGPT is often compared to or personified as a human-like or superintelligent being. This metaphor highlights its sophisticated abilities and the potential for AI to surpass human intelligence in certain areas.' + $vt + '6. **The Future**: The impact of AI on society is frequently discussed in terms of the future, emphasizing that the decisions we make today will shape the world tomorrow. This metaphor acknowledges the long-term implications of AI technology and the importance of careful consideration when developing and deploying it.')

# Question 2
Set-ParagraphText 'The article discusses the role and impact of AI technology, specifically focusing on ChatGPT, in the context of journalism and communication.' ('The topic of AI technology, its impact on society, and its ethical considerations is currently receiving widespread coverage in various media outlets, academic circles, and policy discussions. Some key perspectives that are being widely discussed include:' + $vt + '' + $vt + '1. Advantages of AI technology: Its potential to bring about improvements in healthcare, education, and communication, as well as its ability to facilitate tasks that are difficult or dangerous for humans.' + $vt + '2. Ethical concerns: Issues related to privacy, bias, job displacement, and the need for responsible and ethical use and regulation of AI technology.' + $vt + '3. Future developments: The potential for continued innovation and evolution in AI research and technology, as well as the implications for human employment and societal structures.' + $vt + '' + $vt + 'On the other hand, some aspects that are not being given as much attention include:' + $vt + '' + $vt + '1. Environmental impact: The energy consumption and carbon footprint of AI systems and data centers, and how to minimize their environmental harm.' + $vt + '2. Accessibility and inclusivity: Ensuring that AI technology is accessible and beneficial for people with disabilities, as well as promoting diversity in the development and use of AI.' + $vt + '3. Societal adaptation: The need for education, training, and support for workers who may be displaced by AI automation, as well as developing new models of employment and social safety nets.' + $vt + '4. Cross-cultural considerations: Addressing cultural differences and biases in the development and deployment of AI technology, particularly when working across different geographic regions and cultural contexts.' + $vt + '5. Security and stability: Ensuring that AI systems are secure against hacking and manipulation, as well as mitigating risks associated with AI-driven autonomous weapons and other potentially harmful applications.' + $vt + '6. Legal and regulatory frameworks: Developing clear and effective legal and regulatory frameworks to govern the development, deployment, and use of AI technology, particularly in areas such as healthcare, finance, and public safety.' + $vt + '7. Human-AI collaboration: Exploring new models of human-AI collaboration and partnership, as well as finding ways to promote creativity, innovation, and problem-solving in a world where AI systems are increasingly prevalent.')

# Question 4
Set-ParagraphText 'The article highlights the potential of AI, represented by ChatGPT, in enhancing communication and understanding across cultures' ('The article emphasizes the potential benefits and risks of AI, specifically ChatGPT, in society. It encourages responsible use and ethical regulation to ensure AI technology serves humanity''s best interests.')

# Entities
Set-ParagraphText 'Google, Edge, OpenAI, Arab News, Bard, Microsoft, Meta, AI, Baidu, Bing, Jensen, Amazon, ChatGPT, Frankly Speaking, Katie Jensen, Sam Altman' ('OpenAI, Microsoft, AI, Sam Altman, Katie Jensen, Baidu, Frankly Speaking, Amazon, Google, Bard, Meta, Jensen, Arab News, Bing, ChatGPT, Edge')
